$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '35.490.46'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.60%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.897.37'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.51%  '

# Row 4
$ws.Range('E4').Value = '  +0.02%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.694'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.50%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '246.44'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.20%  '

# Row 7
$ws.Range('E7').Value = '  -0.02%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '43.11'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.97%  '

# Row 9
$ws.Range('E9').Value = '  +1.94%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '56.16'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +8.54%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0758'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.89%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0986'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.48%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.28'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +8.73%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.796'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +9.87%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.172.91'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.54%  '

# Row 16
$ws.Range('E16').Value = '  +2.18%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.923.96'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.43%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '35.455.50'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.66%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '73.74'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.93%  '

# Row 20
$ws.Range('E20').Value = '  +1.45%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '244.64'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.46%  '

# Row 22
$ws.Range('E22').Value = '  +1.51%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.20'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.82%  '

# Row 24
$ws.Range('E24').Value = '  +6.21%  '

# Row 25
$ws.Range('E25').Value = '  +0.05%  '

# Row 26
$ws.Range('E26').Value = '  -1.06%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '166.78'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.93%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.64'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.71%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.33'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.03%  '

# Row 30
$ws.Range('E30').Value = '  +0.62%  '

# Row 31
$ws.Range('E31').Value = '  +2.43%  '

# Row 32
$ws.Range('E32').Value = '  +4.07%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.24'
$ws.Range('D33').Style = 'Normal'

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.88'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +20.58%  '

# Row 35
$ws.Range('E35').Value = '  +0.01%  '

# Row 36
$ws.Range('E36').Value = '  -16.45%  '

# Row 37
$ws.Range('E37').Value = '  +0.04%  '

# Row 38
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0738'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +10.28%  '

# Row 39
$ws.Range('B39').Value = 'LidoDAOToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.95'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.85%  '

# Row 40
$ws.Range('E40').Value = '  +6.27%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '99.22'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.52%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '16.97'
$ws.Range('D42').Style = 'Normal'

# Row 43
$ws.Range('E43').Value = '  -0.98%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.70'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +13.29%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.325.14'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.66%  '

# Row 46
$ws.Range('E46').Value = '  +1.03%  '

# Row 47
$ws.Range('E47').Value = '  +0.36%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.41'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.56%  '

# Row 49
$ws.Range('E49').Value = '  -0.04%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.39'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.27%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '42.59'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.23%  '
